$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price report: insert a new record as row 2 (above the most
# recent existing record), shifting all the existing rows (old 2-18) down
# to become rows 3-19.
$ws.Rows(2).Insert()

# The inserted row picked up the header row's bold/centered formatting;
# strip that so the new row looks like a normal data row again.
$ws.Range("A2:R2").ClearFormats()
# Column D keeps the workbook's date/time number format (same style used
# by every other row in column D).
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45160
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112039
$ws.Range("G2").Value = "Ciboulette"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2500
$ws.Range("N2").Value = "`$/docena de atados"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 833
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = "Hortaliza"
